{"js": "// Replace the 100 addition/subtraction equations in the worksheet table,\n// in row-major order (20 rows x 5 columns), with their new values.\n// Each table cell holds exactly one run with the equation text; we only\n// change the text, not any other formatting.\nconst newValues = [\n  [\"4+84=88\", \"83-43=40\", \"45-23=22\", \"89-81=8\", \"10+31=41\"],\n  [\"15+38=53\", \"69-38=31\", \"88-19=69\", \"1+26=27\", \"94-20=74\"],\n  [\"15+55=70\", \"39+33=72\", \"71+8=79\", \"65+30=95\", \"75-69=6\"],\n  [\"66-53=13\", \"54-38=16\", \"49-23=26\", \"94-77=17\", \"41-34=7\"],\n  [\"89-24=65\", \"19+14=33\", \"52-9=43\", \"35+56=91\", \"34-6=28\"],\n  [\"27-3=24\", \"95+1=96\", \"88-5=83\", \"48-29=19\", \"12+73=85\"],\n  [\"35+16=51\", \"38-30=8\", \"96-16=80\", \"57-17=40\", \"24-7=17\"],\n  [\"86-48=38\", \"49+43=92\", \"32+35=67\", \"79-54=25\", \"62+30=92\"],\n  [\"99-25=74\", \"13+48=61\", \"46-2=44\", \"86-44=42\", \"56-38=18\"],\n  [\"2+84=86\", \"85-8=77\", \"37+52=89\", \"3+10=13\", \"19-10=9\"],\n  [\"34+32=66\", \"72+2=74\", \"44-1=43\", \"28-9=19\", \"65-2=63\"],\n  [\"15+63=78\", \"59-33=26\", \"50-17=33\", \"38+7=45\", \"30-5=25\"],\n  [\"58-17=41\", \"32-2=30\", \"85-84=1\", \"18+40=58\", \"28-22=6\"],\n  [\"35+20=55\", \"23+68=91\", \"23+61=84\", \"8+8=16\", \"86-63=23\"],\n  [\"88-80=8\", \"62+10=72\", \"4+30=34\", \"50+46=96\", \"72-35=37\"],\n  [\"15+74=89\", \"35+52=87\", \"15+7=22\", \"7+26=33\", \"54+18=72\"],\n  [\"52+34=86\", \"53-30=23\", \"4+78=82\", \"52-45=7\", \"53-22=31\"],\n  [\"22+69=91\", \"2+90=92\", \"63-14=49\", \"91-6=85\", \"90-40=50\"],\n  [\"48-35=13\", \"85-33=52\", \"49+27=76\", \"90-81=9\", \"92-62=30\"],\n  [\"17-1=16\", \"24+18=42\", \"60+33=93\", \"37+48=85\", \"18+43=61\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < cells.items.length; c++) {\n    cells.items[c].value = newValues[r][c];\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 addition/subtraction equations in the worksheet table,\n# in row-major order (20 rows x 5 columns), with their new values.\n# Each table cell holds exactly one run with the equation text; only the\n# text is changed, run/paragraph formatting is left untouched.\n$d = $word.ActiveDocument\n$newValues = @(\n    \"4+84=88\", \"83-43=40\", \"45-23=22\", \"89-81=8\", \"10+31=41\",\n    \"15+38=53\", \"69-38=31\", \"88-19=69\", \"1+26=27\", \"94-20=74\",\n    \"15+55=70\", \"39+33=72\", \"71+8=79\", \"65+30=95\", \"75-69=6\",\n    \"66-53=13\", \"54-38=16\", \"49-23=26\", \"94-77=17\", \"41-34=7\",\n    \"89-24=65\", \"19+14=33\", \"52-9=43\", \"35+56=91\", \"34-6=28\",\n    \"27-3=24\", \"95+1=96\", \"88-5=83\", \"48-29=19\", \"12+73=85\",\n    \"35+16=51\", \"38-30=8\", \"96-16=80\", \"57-17=40\", \"24-7=17\",\n    \"86-48=38\", \"49+43=92\", \"32+35=67\", \"79-54=25\", \"62+30=92\",\n    \"99-25=74\", \"13+48=61\", \"46-2=44\", \"86-44=42\", \"56-38=18\",\n    \"2+84=86\", \"85-8=77\", \"37+52=89\", \"3+10=13\", \"19-10=9\",\n    \"34+32=66\", \"72+2=74\", \"44-1=43\", \"28-9=19\", \"65-2=63\",\n    \"15+63=78\", \"59-33=26\", \"50-17=33\", \"38+7=45\", \"30-5=25\",\n    \"58-17=41\", \"32-2=30\", \"85-84=1\", \"18+40=58\", \"28-22=6\",\n    \"35+20=55\", \"23+68=91\", \"23+61=84\", \"8+8=16\", \"86-63=23\",\n    \"88-80=8\", \"62+10=72\", \"4+30=34\", \"50+46=96\", \"72-35=37\",\n    \"15+74=89\", \"35+52=87\", \"15+7=22\", \"7+26=33\", \"54+18=72\",\n    \"52+34=86\", \"53-30=23\", \"4+78=82\", \"52-45=7\", \"53-22=31\",\n    \"22+69=91\", \"2+90=92\", \"63-14=49\", \"91-6=85\", \"90-40=50\",\n    \"48-35=13\", \"85-33=52\", \"49+27=76\", \"90-81=9\", \"92-62=30\",\n    \"17-1=16\", \"24+18=42\", \"60+33=93\", \"37+48=85\", \"18+43=61\"\n)\n\n$t = $d.Tables.Item(1)\n$cols = $t.Columns.Count\n$idx = 0\nforeach ($row in $t.Rows) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($row.Index, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
